# "emission plots and postprocess" -- emission default name of CO2
# equivalent has changed to a user-specified emission name.
#
# On the "Emissions" sheet we insert a new "emission_name" column between
# the emission id column (A) and the "emission_unit" column (shifted from
# B to C), rename the emission id value from "CO2-equivalent" to "CO2",
# and give it a separate descriptive name "CO2 emissions" in the new
# column. The "Regions" sheet loses the tab-selection (it moves to
# "Emissions" instead), and column widths are set to fit the new data.

$wb = $excel.ActiveWorkbook
$wsReg = $wb.Worksheets.Item("Regions")
$ws = $wb.Worksheets.Item("Emissions")

# Make room for the new column: shift emission_unit/ton from column B to
# column C (preserves their existing formatting).
$ws.Range("B1").EntireColumn.Insert()

# Give the new header cell (B1) and the shifted unit header (C1) the same
# styling as the rest of the header row by copying an already-styled
# header cell onto them, then overwrite with the new text.
$ws.Range("C1").Copy($ws.Range("B1"))
$ws.Range("B1").Value = "emission_name"
$ws.Range("C1").Value = "emission_unit"

# A2 keeps its existing (styled) position; just change its text.
$ws.Range("A2").Value = "CO2"

# B2 (new "CO2 emissions" descriptive name) should be unstyled, like C2 --
# copy C2's (unstyled) formatting onto it before setting its text.
$ws.Range("C2").Copy($ws.Range("B2"))
$ws.Range("B2").Value = "CO2 emissions"
$ws.Range("C2").Value = "ton"

# Column widths to fit the new content.
$ws.Columns.Item(2).ColumnWidth = 13.1640625
$ws.Columns.Item(3).ColumnWidth = 12

# The "Emissions" tab becomes the selected/active one (was "Regions").
$wsReg.Range("A1").Select()
$ws.Range("E6").Select()
$ws.Activate()
